$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column to a Text format before writing, so that
# numeric-looking strings (e.g. "593.66", "27.00", "3.281.06") are stored
# verbatim as text instead of being re-parsed/rounded as numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '67.477.07'
$ws.Range("E2").Value = '  -4.62%  '

# Row 3
$ws.Range("D3").Value = '3.281.06'
$ws.Range("E3").Value = '  -7.11%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '593.66'
$ws.Range("E5").Value = '  -4.36%  '

# Row 6
$ws.Range("D6").Value = '152.69'
$ws.Range("E6").Value = '  -11.61%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").Value = '3.272.82'
$ws.Range("E8").Value = '  -7.22%  '

# Row 9
$ws.Range("E9").Value = '  -10.45%  '

# Row 10
$ws.Range("E10").Value = '  -13.28%  '

# Row 11
$ws.Range("D11").Value = '6.75'
$ws.Range("E11").Value = '  -6.46%  '

# Row 12
$ws.Range("D12").Value = '0.513'
$ws.Range("E12").Value = '  -12.52%  '

# Row 13
$ws.Range("D13").Value = '38.84'
$ws.Range("E13").Value = '  -15.98%  '

# Row 14
$ws.Range("E14").Value = '  -10.80%  '

# Row 15
$ws.Range("D15").Value = '3.810.06'
$ws.Range("E15").Value = '  -7.07%  '

# Row 16
$ws.Range("D16").Value = '67.497.67'
$ws.Range("E16").Value = '  -4.73%  '

# Row 17
$ws.Range("D17").Value = '3.284.71'
$ws.Range("E17").Value = '  -7.15%  '

# Row 18
$ws.Range("E18").Value = '  -13.67%  '

# Row 19
$ws.Range("D19").Value = '536.83'
$ws.Range("E19").Value = '  -11.62%  '

# Row 21
$ws.Range("D21").Value = '15.15'
$ws.Range("E21").Value = '  -14.35%  '

# Row 22
$ws.Range("E22").Value = '  -13.16%  '

# Row 23
$ws.Range("D23").Value = '7.92'
$ws.Range("E23").Value = '  -12.75%  '

# Row 24
$ws.Range("D24").Value = '13.73'
$ws.Range("E24").Value = '  -12.14%  '

# Row 25
$ws.Range("D25").Value = '86.02'
$ws.Range("E25").Value = '  -11.80%  '

# Row 26
$ws.Range("E26").Value = '  -0.05%  '

# Row 27
$ws.Range("E27").Value = '  -11.74%  '

# Row 28
$ws.Range("D28").Value = '8.19'
$ws.Range("E28").Value = '  -10.13%  '

# Row 29
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").Value = '2.17'
$ws.Range("E29").Value = '  -15.49%  '

# Row 30
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '29.58'
$ws.Range("E30").Value = '  -12.08%  '

# Row 31
$ws.Range("E31").Value = '  -10.39%  '

# Row 32
$ws.Range("D32").Value = '1.16'
$ws.Range("E32").Value = '  -10.61%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '6.65'
$ws.Range("E33").Value = '  -18.00%  '

# Row 34
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").Value = '542.78'
$ws.Range("E34").Value = '  -12.09%  '

# Row 35
$ws.Range("D35").Value = '5.84'
$ws.Range("E35").Value = '  -14.37%  '

# Row 36
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.21%  '

# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.0467'
$ws.Range("E37").Value = '  -5.92%  '

# Row 38
$ws.Range("D38").Value = '53.54'
$ws.Range("E38").Value = '  -5.71%  '

# Row 39
$ws.Range("D39").Value = '0.0866'
$ws.Range("E39").Value = '  -12.90%  '

# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '0.130'
$ws.Range("E40").Value = '  -9.45%  '

# Row 41
$ws.Range("B41").Value = 'Cosmos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D41").Value = '9.16'
$ws.Range("E41").Value = '  -15.49%  '

# Row 42
$ws.Range("D42").Value = '2.76'
$ws.Range("E42").Value = '  -18.47%  '

# Row 43
$ws.Range("D43").Value = '2.952.90'
$ws.Range("E43").Value = '  -11.85%  '

# Row 44
$ws.Range("E44").Value = '  -12.48%  '

# Row 45
$ws.Range("D45").Value = '0.0₃0599'
$ws.Range("E45").Value = '  -17.44%  '

# Row 46
$ws.Range("D46").Value = '2.23'
$ws.Range("E46").Value = '  -10.85%  '

# Row 47
$ws.Range("D47").Value = '27.00'
$ws.Range("E47").Value = '  -15.21%  '

# Row 48
$ws.Range("E48").Value = '  -17.83%  '

# Row 49
$ws.Range("E49").Value = '  -0.11%  '

# Row 50
$ws.Range("D50").Value = '125.79'
$ws.Range("E50").Value = '  -6.15%  '

# Row 51
$ws.Range("E51").Value = '  -11.84%  '

# Restore the default (unstyled) cell style now that the text values are
# committed, so the cells end up styled exactly as they were originally.
$priceRange.Style = "Normal"
